# Update the "dSF" column (F) values on Sheet1 with freshly re-pulled data.
# Mirrors the underlying data refresh described in the commit message
# ("repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F2"  = 1
    "F3"  = -3
    "F4"  = -1
    "F5"  = 2
    "F6"  = -1
    "F8"  = 6
    "F9"  = 1
    "F10" = -3
    "F11" = -3
    "F12" = -4
    "F13" = -3
    "F15" = -5
    "F16" = 4
    "F17" = -5
    "F18" = -4
    "F19" = 4
    "F20" = 2
    "F21" = 5
    "F22" = -1
    "F23" = -3
    "F24" = 1
    "F26" = 6
    "F28" = 1
    "F29" = -1
    "F30" = 1
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
